$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.526.81'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.60%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.442.29'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -4.91%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.18%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.44'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -5.08%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '189.27'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.73%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.604'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.48%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.432.55'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -4.83%  '

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -5.24%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.612'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -5.25%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '50.64'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.89%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000282'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -7.29%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.03'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -5.58%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.990.90'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -4.86%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '631.45'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +4.14%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.376.47'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.90%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.443.31'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -4.41%  '

$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.120'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.54%  '

$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.18'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -5.60%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.98'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -5.66%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.934'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -6.63%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.70'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.65%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.34'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.26%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '98.55'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -4.33%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.23'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -8.29%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.81'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -6.21%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.06'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.74%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.78'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -8.34%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.11'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -5.82%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.07'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -5.09%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.12'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -11.83%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.66'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -8.34%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.49'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -6.59%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '60.62'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.28%  '

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -7.71%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.06%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.629.28'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -7.40%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0₃0777'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -12.29%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '500.07'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.04%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.46'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.19%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.88'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -6.25%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.364'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -6.27%  '

$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.132'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.17%  '

$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '34.06'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -7.51%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.43'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +61.23%  '

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -5.39%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.34'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.29%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.78'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -4.06%  '

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.78%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.26%  '
